$wb = $excel.ActiveWorkbook

# Property sheet: rename "View" field to "Cache"
$wsProperty = $wb.Worksheets.Item("Property")
$wsProperty.Range("F1").Value = "Cache"

# Record sheet shares the same "View" header -> rename it too
$wsRecord = $wb.Worksheets.Item("Record")
$wsRecord.Range("G1").Value = "Cache"

# Set all of F2:F26 to FALSE (was TRUE)
$wsProperty.Range("F2:F26").Value = $false

# F24:F26 previously carried a highlight fill (style) that should no longer
# apply once the whole column is uniformly FALSE; re-stamp the format of the
# (unstyled) F2 cell across the whole column so it matches.
$wsProperty.Range("F2").Copy($wsProperty.Range("F2:F26"))

# Select F2:F26 on Property sheet, with F2 as the active cell
$wsProperty.Range("F2:F26").Select()

# Make Property the active sheet/tab
$wsProperty.Activate()
